# parser fetch total time info
#
# - Update the instructional header text (C1) to reference "a node id"
#   instead of "an X".
# - Add a new "t02n41" node-id value in C5, and extend the node-id
#   column (C:F) with the same right-aligned / General-format styling
#   down through row 11, for every data row (5-11).
# - Move the active selection to C2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header / instructions cell text.
$ws.Range("C1").Value = "if there’s a node id, then we ran it!"

# New node-id column: right aligned, General number format (matches the
# existing "8-1-1"/"4-1-1" style but without the forced text format).
$ws.Range("C5:F11").HorizontalAlignment = -4152
$ws.Range("C5:F11").NumberFormat = "General"

# Only the first row of the new column actually carries a value.
$ws.Range("C5").Value = "t02n41"

# Move the selection cursor to C2.
$ws.Range("C2").Select()
